$wb = $excel.ActiveWorkbook
$texts = $wb.Worksheets.Item("texts")
$fc = $wb.Worksheets.Item("functionsCategories")

# ---------------------------------------------------------------------------
# 1) texts sheet: tighten the drop-shadow filter (dx/dy/stdDeviation 5 -> 3)
# ---------------------------------------------------------------------------
$svgBase = $texts.Range("A1").Value2
$svgBase = $svgBase -replace "dx='5' dy='5'", "dx='3' dy='3'"
$svgBase = $svgBase -replace "stdDeviation='5'", "stdDeviation='3'"
$texts.Range("A1").Value = $svgBase

# ---------------------------------------------------------------------------
# 2) functionsCategories sheet: refresh the per-category fill/stroke colors
#    (column E = fill color, column F = stroke color)
# ---------------------------------------------------------------------------
$fc.Range("E2").Value = "#FF0000"
$fc.Range("F2").Value = "#FF3333"

$fc.Range("E3").Value = "#80FF00"
$fc.Range("F3").Value = "#99FF33"

$fc.Range("E4").Value = "#00FFFF"
$fc.Range("F4").Value = "#33FFFF"

$fc.Range("E5").Value = "#8000FF"
$fc.Range("F5").Value = "#9933FF"

$fc.Range("E6").Value = "#FFFF00"
$fc.Range("F6").Value = "#FFFF33"

$fc.Range("E7").Value = "#00FF80"
$fc.Range("F7").Value = "#33FF99"

$fc.Range("E8").Value = "#0000FF"
$fc.Range("F8").Value = "#3333FF"

$fc.Range("E9").Value = "#FF007F"
$fc.Range("F9").Value = "#FF3399"

# ---------------------------------------------------------------------------
# 3) functionsCategories sheet: new preIcon / postIcon helper columns (L, M)
# ---------------------------------------------------------------------------
$fc.Range("L1").Value = "preIcon"
$fc.Range("M1").Value = "postIcon"

for ($r = 2; $r -le 8; $r++) {
    $fc.Range("L$r").Formula = "=texts!`$A`$1&functionsCategories!E$r&texts!`$A`$2&functionsCategories!F$r&texts!`$A`$3"
    $fc.Range("M$r").Formula = "=texts!`$A`$4"
}

# ---------------------------------------------------------------------------
# 4) Cell formatting: give the texts!A2 cell the same quote-prefix style as
#    the other formula-driven string cells (keeps the same wrap formatting)
# ---------------------------------------------------------------------------
$texts.Range("A2").NumberFormat = "@"
$texts.Range("A2").WrapText = $true

# ---------------------------------------------------------------------------
# 5) Selection / active-sheet bookkeeping
# ---------------------------------------------------------------------------
$texts.Activate()
$texts.Range("A3").Select()

$fc.Activate()
$fc.Range("L2").Select()
